# Applies updated per-bus voltage magnitude (vm_pu) results for the
# "380 kV" case (Case_1_74) onto the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..F (bus indices 0-4) for rows 2..25 (time steps 0..23)
$block1Rows = @(
    ,(1.02, 1.036487756214196, 1.057963630632712, 1.035389194267288, 1.063173931094379)
    ,(1.019999999999999, 1.037528012924204, 1.058701374171659, 1.036276702230254, 1.064069889248569)
    ,(1.02, 1.038201469070043, 1.059178994806896, 1.036851610875268, 1.064650230802678)
    ,(1.02, 1.038484671321531, 1.059379845415652, 1.037093452881513, 1.064894348042456)
    ,(1.02, 1.038532226985992, 1.059413572533787, 1.037134068015367, 1.064935344672626)
    ,(1.02, 1.038205252913251, 1.059181678352719, 1.036854841789981, 1.064653492154102)
    ,(1.02, 1.036839245200082, 1.058212901743566, 1.035689000467232, 1.06347660006865)
    ,(1.019999999999999, 1.034434782687954, 1.056507772929563, 1.03363951219015, 1.061407389098932)
    ,(1.02, 1.032833578403082, 1.055372424791588, 1.032276509348269, 1.06003109402473)
    ,(1.02, 1.032140660725594, 1.054881154880085, 1.03168711144629, 1.059435911564501)
    ,(1.02, 1.031883342515138, 1.054698728140759, 1.031468302403063, 1.059214949967215)
    ,(1.02, 1.031938535326722, 1.054737856887979, 1.031515232263186, 1.059262341727516)
    ,(1.02, 1.032119389451697, 1.054866074345285, 1.031669022179097, 1.059417644443799)
    ,(1.02, 1.032230827898639, 1.054945080382097, 1.03176379308834, 1.059513346933964)
    ,(1.02, 1.032879573798301, 1.055405036093475, 1.032315642477289, 1.060070610483416)
    ,(1.02, 1.03328662590323, 1.055693646923845, 1.032662015367918, 1.06042037209512)
    ,(1.02, 1.033524092486047, 1.055862021783778, 1.032864125119501, 1.060624455609519)
    ,(1.02, 1.033605069234437, 1.055919438828412, 1.032933052243328, 1.060694055243733)
    ,(1.02, 1.033242948929473, 1.055662678286394, 1.032624844937185, 1.060382838375333)
    ,(1.02, 1.032066130700522, 1.054828316033245, 1.031623731575393, 1.059371908462577)
    ,(1.02, 1.031326578643292, 1.054304025328224, 1.030994984071013, 1.058736966486668)
    ,(1.02, 1.031718595024979, 1.054581932349383, 1.03132822917495, 1.059073497376835)
    ,(1.02, 1.033262684550216, 1.055676671578635, 1.032641640419097, 1.06039979802072)
    ,(1.02, 1.035056082589305, 1.056948347448377, 1.034168771784538, 1.061941774591731)
)

# Columns I..N (bus indices 7-12) for rows 2..25 (time steps 0..23)
$block2Rows = @(
    ,(1.044836911379219, 1.041595802792134, 1.060697162094288, 1.038185939347275, 1.065893285383269, 1.017750063593984)
    ,(1.04510500594365, 1.042279587321826, 1.061249038772682, 1.038882624337137, 1.066603999489315, 1.017982272930095)
    ,(1.045277487504928, 1.042721805842332, 1.061605721642464, 1.039333423842552, 1.067063837685844, 1.018132317297507)
    ,(1.045349760465409, 1.042907657536701, 1.061755570095419, 1.0395229388102, 1.067257142777589, 1.018195345249794)
    ,(1.045361881427711, 1.042938859498008, 1.061780724359791, 1.039554759122635, 1.067289598886079, 1.018205924939082)
    ,(1.045278454156604, 1.042724289426276, 1.061607724322909, 1.039335956156072, 1.067066420681948, 1.018133159679657)
    ,(1.044927720654147, 1.041826939483495, 1.060883757462425, 1.038421387409397, 1.066133482512753, 1.017828583263229)
    ,(1.044302092996939, 1.040243911253768, 1.059604871641568, 1.036809809458926, 1.064489252975374, 1.017290280124944)
    ,(1.043879929523052, 1.039187389762967, 1.058750206481247, 1.03573546448051, 1.063392973742324, 1.016930351427096)
    ,(1.043695928616672, 1.038729632664543, 1.058379645757979, 1.03527027758277, 1.062918254589901, 1.016774250181132)
    ,(1.043627402368699, 1.038559560178605, 1.058241930960166, 1.035097488682868, 1.062741920208484, 1.01671622994693)
    ,(1.043642109626319, 1.038596043158449, 1.058271474498678, 1.035134552379832, 1.062779744632782, 1.016728677167307)
    ,(1.043690267887865, 1.038715575256293, 1.058368263674743, 1.035255994750152, 1.06290367878173, 1.016769454969605)
    ,(1.043719915925814, 1.038789217490204, 1.058427889158616, 1.035330819712273, 1.062980038390271, 1.016794574586529)
    ,(1.043892115777842, 1.039217763805991, 1.058774789227039, 1.035766337688222, 1.06342447888151, 1.016940706111355)
    ,(1.04399981082318, 1.039486505887805, 1.058992261194146, 1.036039530049572, 1.063703259138363, 1.017032303861991)
    ,(1.0440625115847, 1.039643231909705, 1.059119062055409, 1.036198879615139, 1.063865864664671, 1.017085707167374)
    ,(1.044083871248048, 1.039696666894304, 1.059162289899126, 1.036253213867939, 1.063921308532612, 1.017103912217178)
    ,(1.043988268153192, 1.03945767516329, 1.05896893335091, 1.036010218979041, 1.063673348884502, 1.017022478782573)
    ,(1.04367609144934, 1.038680377164599, 1.05833976365695, 1.035220232933573, 1.062867183326993, 1.016757447950398)
    ,(1.043478771638198, 1.038191421364322, 1.057943763274343, 1.034723550476704, 1.062360300732178, 1.016590596969644)
    ,(1.043583473242517, 1.038450648483416, 1.058153729764593, 1.03498684983698, 1.062629009840201, 1.016679068185451)
    ,(1.043993484147146, 1.039470702609127, 1.058979474343977, 1.036023463386777, 1.063686864047858, 1.017026918387557)
    ,(1.044464729677444, 1.040653370263245, 1.059935863790368, 1.037226436394743, 1.064914352182591, 1.017429632397003)
)

$numRows = $block1Rows.Count
$block1 = New-Object "object[,]" $numRows,5
$block2 = New-Object "object[,]" $numRows,6
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt 5; $j++) { $block1[$i,$j] = $block1Rows[$i][$j] }
    for ($j = 0; $j -lt 6; $j++) { $block2[$i,$j] = $block2Rows[$i][$j] }
}

$ws.Range("B2:F25").Value = $block1
$ws.Range("I2:N25").Value = $block2

Write-Output "Applied vm_pu updates to B2:F25 and I2:N25"
